$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 is "CETAL 500MG 20 TAB" - a sale was recorded against it:
#   - current balance (H9) moved from 0:1 to 0:0
#   - selling price (P9) doubled from 12.0000 to 24.0000
#   - transaction count (Q9) flipped from 0:1 to 1:0
$ws.Range("H9").Value = "0:0"
$ws.Range("Q9").Value = "1:0"

# P9 is formatted as a number (0.00) but the sheet stores it as text, so
# round-trip the number format through "@" to keep it a literal string
# "24.0000" instead of letting COM coerce it into the numeric 24.
$fmt = $ws.Range("P9").NumberFormat
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "24.0000"
$ws.Range("P9").NumberFormat = $fmt

# Grand total (P17) increased to reflect the new selling price
$ws.Range("P17").Value = 443.73000000000002

# Footer timestamp bumped by two minutes on re-export
$ws.Range("A18").Value = "Saturday, 31 May, 2025 10:03 AM"
